$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells stay formatted as text (matching the original inline-string
# representation) so the values are not auto-converted to numbers/percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "274.76"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.15%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.37%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.872"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.27%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.34%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.861"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.88%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.311"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.40%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.246"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "32.25%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8691"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.23%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1684"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "16.07%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05038"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.33%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07423"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.06%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02964"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.57%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09025"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.34%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001580"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.43%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006312"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.54%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005933"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.29%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.446"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.17%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.34%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.25%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1334"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.48%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.918"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.74%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04356"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.74%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001177"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.15%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.90%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.01%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001687"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.22%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04046"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.08%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006720"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.80%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1166"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.94%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002191"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.63%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-11.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005301"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.90%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.02101"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-29.67%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.486"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-37.43%"
